$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data block (rows 2-9, columns A-H) down to rows 10-17
# — this duplicates the roster so it repeats for a second group, same as
# the source workbook's data range A2:H9 -> A10:H17.
$src = $ws.Range("A2:H9")
$dst = $ws.Range("A10:H17")
$src.Copy($dst)

# Fix up the "NO" column so numbering continues (9..16) instead of repeating 1..8
for ($i = 0; $i -lt 8; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = 9 + $i
}

# Match the reported selection after the edit
$ws.Range("J7").Select()
